$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to advance by one
# day (46061 -> 46062) for every data row (rows 2 through 75).
for ($row = 2; $row -le 75; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 46062
}
